$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53; this shifts existing rows 53-70 down to 54-71
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new data record
$ws.Range("A53").Value = 4
$ws.Range("B53").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C53").Value = "Los Lagos"
$ws.Range("D53").Value = 44463
$ws.Range("E53").Value = 10
$ws.Range("F53").Value = 100112022
$ws.Range("G53").Value = "Arveja Verde"
$ws.Range("H53").Value = "Perfection"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 70
$ws.Range("K53").Value = 35000
$ws.Range("L53").Value = 35000
$ws.Range("M53").Value = 35000
$ws.Range("N53").Value = "$/malla 25 kilos"
$ws.Range("O53").Value = "Provincia de Huasco"
$ws.Range("P53").Value = 1400
$ws.Range("Q53").Value = 25
$ws.Range("R53").Value = "Hortaliza"
